# feat: add 2022-Q1 data
#
# 1) Duplicate the "2021-Q4" sheet (it already has the right column layout /
#    styles for a per-fund holdings sheet) and drop the copy right after it
#    -- i.e. right before "总计" -- then rename it to "2022-Q1" and overwrite
#    its numbers with the new quarter's figures.
# 2) Insert a new first data row into "总计" for 2022-Q1 and renumber the
#    existing rows' index column.

$wb = $excel.ActiveWorkbook

# --- 1) New "2022-Q1" sheet ------------------------------------------------

$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy([System.Reflection.Missing]::Value, $template)
$ws = $template.Next
$ws.Name = "2022-Q1"

# Fund rows: code, name, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$funds = @(
    @("160613", "鹏华盛世创新混合(LOF)",                 "3.05", "93.82", "4.52", "0.1379", 7),
    @("011574", "鹏华领航一年持有期混合型证券投资基金A", "2.75", "90.27", "4.03", "0.1108", 7),
    @("008134", "鹏华优选价值股票",                       "1.86", "92.62", "4.20", "0.0781", 6),
    @("011575", "鹏华领航一年持有期混合型证券投资基金C", "1.22", "90.27", "4.03", "0.0492", 7)
)

for ($i = 0; $i -lt $funds.Count; $i++) {
    $row = $i + 2
    $f = $funds[$i]
    $ws.Cells.Item($row, 2).Value = "'" + $f[0]
    $ws.Cells.Item($row, 3).Value = $f[1]
    $ws.Cells.Item($row, 4).Value = "'" + $f[2]
    $ws.Cells.Item($row, 5).Value = "'" + $f[3]
    $ws.Cells.Item($row, 6).Value = "'" + $f[4]
    $ws.Cells.Item($row, 7).Value = "'" + $f[5]
    $ws.Cells.Item($row, 8).Value = $f[6]
}

# --- 2) Update the "总计" sheet ---------------------------------------------

$tot = $wb.Worksheets.Item("总计")
$tot.Rows.Item(2).Insert()

# Row-insert pulls its formatting from the header row above (bold, no
# border) -- strip that and instead clone the plain "index column" look
# from the data row right below (same as A3:A7).
$tot.Rows.Item(2).ClearFormats()
$tot.Cells.Item(3, 1).Copy()
$tot.Cells.Item(2, 1).PasteSpecial(-4122)  # xlPasteFormats

# New first row: 2022-Q1
$tot.Cells.Item(2, 1).Value = 0
$tot.Cells.Item(2, 2).Value = "2022-Q1"
$tot.Cells.Item(2, 3).Value = 4
$tot.Cells.Item(2, 4).Value = 0.38

# Renumber the index column for the pre-existing rows (now shifted down one).
for ($r = 3; $r -le 7; $r++) {
    $tot.Cells.Item($r, 1).Value = $r - 2
}

# Restore the originally-active sheet (copying a sheet makes the new copy
# active).
$wb.Worksheets.Item(1).Activate()
